$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: numeric value refresh
$ws.Range("B2").Value = 0.5608

# C2: text-typed numeric string refresh (force text storage, then drop the
# transient "Text" number format style so no stray style id is left on the cell)
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "0.9147"
$ws.Range("C2").Style = "Normal"

# D2: text-typed numeric string refresh
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "0.5608"
$ws.Range("D2").Style = "Normal"

# E2: numeric value refresh
$ws.Range("E2").Value = 0.6860000000000001
